# Update "想去人数" (attendance count) values on both the "展览" sheet and
# the combined "全部类型" sheet to reflect newly scraped figures.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14838
$ws1.Range("F3").Value = 18405
$ws1.Range("F5").Value = 110
$ws1.Range("F14").Value = 97
$ws1.Range("F15").Value = 194
$ws1.Range("F17").Value = 1402
$ws1.Range("F20").Value = 82
$ws1.Range("F21").Value = 226
$ws1.Range("F22").Value = 7624
$ws1.Range("F24").Value = 19
$ws1.Range("F25").Value = 50
$ws1.Range("F26").Value = 1210
$ws1.Range("F28").Value = 5940
$ws1.Range("F29").Value = 94
$ws1.Range("F30").Value = 60
$ws1.Range("F34").Value = 5275
$ws1.Range("F35").Value = 24
$ws1.Range("F36").Value = 38

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14838
$ws4.Range("F3").Value = 18405
$ws4.Range("F5").Value = 110
$ws4.Range("F14").Value = 97
$ws4.Range("F15").Value = 194
$ws4.Range("F17").Value = 1402
$ws4.Range("F21").Value = 82
$ws4.Range("F22").Value = 226
$ws4.Range("F23").Value = 7624
$ws4.Range("F25").Value = 19
$ws4.Range("F26").Value = 50
$ws4.Range("F27").Value = 1210
$ws4.Range("F30").Value = 5940
$ws4.Range("F31").Value = 94
$ws4.Range("F32").Value = 60
$ws4.Range("F36").Value = 5275
$ws4.Range("F37").Value = 24
$ws4.Range("F38").Value = 38
